$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column F header "Static"
$ws.Range("F1").Value = "Static"

# Fill F2:F26 with 0 (Wh -> kWh conversion factor column, static placeholder values)
for ($r = 2; $r -le 26; $r++) {
    $ws.Cells.Item($r, 6).Value = 0
}

# Update the selected cell to match the author's saved selection
$ws.Range("C31").Select()
